$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 300
$ws.Range("F2").Value = 60

$ws.Range("E2").Select()
